$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999999999915534
$ws.Range("E2").Value = 0.9999999999915534

# Row 3
$ws.Range("D3").Value = 0.9999998040056253
$ws.Range("E3").Value = 0.9999998040056253

# Row 4
$ws.Range("D4").Value = 0.370677866201541
$ws.Range("E4").Value = 0.370677866201541

# Row 6
$ws.Range("D6").Value = 0.9970115336646616
$ws.Range("E6").Value = 0.9970115336646616

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.2484036296847743
$ws.Range("E7").Value = 0.7515963703152256

# Row 8
$ws.Range("D8").Value = 0.9999999966082147
$ws.Range("E8").Value = [double]"3.39178529706885E-09"

# Row 9
$ws.Range("D9").Value = 0.9999994615271761
$ws.Range("E9").Value = [double]"5.384728238855274E-07"

# Row 10
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0

# Row 11
$ws.Range("D11").Value = 0.9999999999999951
$ws.Range("E11").Value = [double]"4.884981308350689E-15"
$ws.Range("F11").Value = 9.469118118286133
$ws.Range("G11").Value = 0.5

# Row 13
$ws.Range("D13").Value = 0.9999999867922651
$ws.Range("E13").Value = 0.9999999867922651

# Row 14
$ws.Range("D14").Value = 0.3920504207009849
$ws.Range("E14").Value = 0.3920504207009849

# Row 16
$ws.Range("D16").Value = 0.9997111999204039
$ws.Range("E16").Value = 0.9997111999204039

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = 0.01493643850081175
$ws.Range("E17").Value = 0.9850635614991883

# Row 18
$ws.Range("D18").Value = 0.9999999999973488
$ws.Range("E18").Value = [double]"2.651212582804874E-12"

# Row 19
$ws.Range("D19").Value = 0.9999999995463422
$ws.Range("E19").Value = [double]"4.53657778010097E-10"

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0

# Row 21
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 12.90282440185547
$ws.Range("G21").Value = 0.5
